# Update the cryptos price/volume table with freshly scraped values.
# Row 11-13 also get re-sorted (Solana moves up, TRON moves up, WrappedEther moves down),
# per the source coinranking.com ordering for this run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Cell holds a decimal-looking string (e.g. "1.007") that must stay TEXT,
    # not get auto-coerced to a number by Excel. A leading apostrophe forces
    # text entry; we then reset the cell style back to Normal so no stray
    # "quote prefix" number format lingers on the cell.
    $ws.Range($range).Value = "'" + $value
    $ws.Range($range).Style = "Normal"
}

Set-TextValue "D2" "26.887.22"
$ws.Range("E2").Value = "  +1.26%  "

Set-TextValue "D3" "1.841.53"
$ws.Range("E3").Value = "  +1.61%  "

Set-TextValue "D5" "309.48"
$ws.Range("E5").Value = "  +1.16%  "

Set-TextValue "D6" "1.007"
$ws.Range("E6").Value = "  +0.39%  "

Set-TextValue "D7" "0.4715"
$ws.Range("E7").Value = "  +3.62%  "

Set-TextValue "D8" "0.3650"
$ws.Range("E8").Value = "  +1.49%  "

Set-TextValue "D9" "0.07155"
$ws.Range("E9").Value = "  +0.62%  "

Set-TextValue "D10" "0.9198"
$ws.Range("E10").Value = "  +2.97%  "

# Rows 11-13 swap identities: WrappedEther/Solana/TRON -> Solana/TRON/WrappedEther
$ws.Range("B11").Value = "Solana"
$ws.Range("C11").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
Set-TextValue "D11" "19.52"
$ws.Range("E11").Value = "  +0.97%  "

$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D12" "0.07608"
$ws.Range("E12").Value = "  -1.45%  "

$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
Set-TextValue "D13" "1.852.50"
$ws.Range("E13").Value = "  +2.43%  "

Set-TextValue "D14" "5.281"
$ws.Range("E14").Value = "  +0.51%  "

Set-TextValue "D15" "6.393"
$ws.Range("E15").Value = "  +1.62%  "

Set-TextValue "D16" "87.83"
$ws.Range("E16").Value = "  +1.33%  "

Set-TextValue "D17" "1.009"
$ws.Range("E17").Value = "  +0.39%  "

Set-TextValue "D18" "0.000008628"
$ws.Range("E18").Value = "  +0.91%  "

Set-TextValue "D19" "1.007"
$ws.Range("E19").Value = "  +0.38%  "

Set-TextValue "D20" "26.907.88"

Set-TextValue "D21" "14.48"
$ws.Range("E21").Value = "  +2.34%  "

Set-TextValue "D22" "5.012"
$ws.Range("E22").Value = "  +1.03%  "

Set-TextValue "D23" "10.61"
$ws.Range("E23").Value = "  +0.74%  "

Set-TextValue "D24" "1.929"
$ws.Range("E24").Value = "  +0.43%  "

Set-TextValue "D25" "151.62"
$ws.Range("E25").Value = "  -0.06%  "

Set-TextValue "D26" "18.19"
$ws.Range("E26").Value = "  +2.14%  "

Set-TextValue "D27" "2.009"
$ws.Range("E27").Value = "  -0.43%  "

Set-TextValue "D28" "114.12"
$ws.Range("E28").Value = "  +1.46%  "

Set-TextValue "D29" "4.856"
$ws.Range("E29").Value = "  +0.50%  "

Set-TextValue "D30" "0.08821"

Set-TextValue "D31" "3.230"
$ws.Range("E31").Value = "  +3.40%  "

Set-TextValue "D32" "1.169"
$ws.Range("E32").Value = "  +5.27%  "

$ws.Range("E33").Value = "  +0.39%  "

Set-TextValue "D34" "4.475"
$ws.Range("E34").Value = "  +0.89%  "

Set-TextValue "D35" "2.749"
$ws.Range("E35").Value = "  +0.73%  "

Set-TextValue "D36" "1.089"
$ws.Range("E36").Value = "  +1.77%  "

Set-TextValue "D37" "0.01942"
$ws.Range("E37").Value = "  +0.32%  "

Set-TextValue "D38" "0.05232"
$ws.Range("E38").Value = "  +3.12%  "

Set-TextValue "D39" "2.967"
$ws.Range("E39").Value = "  +1.79%  "

Set-TextValue "D40" "0.5178"
$ws.Range("E40").Value = "  +1.89%  "

Set-TextValue "D41" "6.949"
$ws.Range("E41").Value = "  +2.49%  "

Set-TextValue "D42" "0.1511"
$ws.Range("E42").Value = "  +0.27%  "

Set-TextValue "D43" "8.148"
$ws.Range("E43").Value = "  +1.64%  "

Set-TextValue "D44" "10.47"
$ws.Range("E44").Value = "  +4.99%  "

Set-TextValue "D45" "0.4699"
$ws.Range("E45").Value = "  +0.25%  "

Set-TextValue "D46" "1.008"
$ws.Range("E46").Value = "  +0.43%  "

Set-TextValue "D47" "102.05"
$ws.Range("E47").Value = "  +2.91%  "

Set-TextValue "D48" "1.594"

Set-TextValue "D49" "65.56"
$ws.Range("E49").Value = "  +3.08%  "

Set-TextValue "D50" "0.06035"
$ws.Range("E50").Value = "  +0.55%  "

Set-TextValue "D51" "0.8858"
$ws.Range("E51").Value = "  +4.55%  "
